# Update the multiplication-fact answers in the single table of the
# document. Each data row (1, 5, 10, 15, 20 in 1-based Word Table.Cell
# indexing) holds five "A×B=C" strings across its five columns; replace
# them in place with the new facts, preserving all run/paragraph
# formatting by writing only into Cell.Range.Text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "337×4=1348" },
    @{ Row = 1;  Col = 2; New = "991×9=8919" },
    @{ Row = 1;  Col = 3; New = "430×9=3870" },
    @{ Row = 1;  Col = 4; New = "378×5=1890" },
    @{ Row = 1;  Col = 5; New = "946×8=7568" },

    @{ Row = 5;  Col = 1; New = "893×8=7144" },
    @{ Row = 5;  Col = 2; New = "691×8=5528" },
    @{ Row = 5;  Col = 3; New = "409×7=2863" },
    @{ Row = 5;  Col = 4; New = "726×4=2904" },
    @{ Row = 5;  Col = 5; New = "825×4=3300" },

    @{ Row = 10; Col = 1; New = "301×6=1806" },
    @{ Row = 10; Col = 2; New = "170×8=1360" },
    @{ Row = 10; Col = 3; New = "894×9=8046" },
    @{ Row = 10; Col = 4; New = "980×4=3920" },
    @{ Row = 10; Col = 5; New = "525×9=4725" },

    @{ Row = 15; Col = 1; New = "435×8=3480" },
    @{ Row = 15; Col = 2; New = "284×4=1136" },
    @{ Row = 15; Col = 3; New = "173×3=519"  },
    @{ Row = 15; Col = 4; New = "147×4=588"  },
    @{ Row = 15; Col = 5; New = "449×5=2245" },

    @{ Row = 20; Col = 1; New = "387×7=2709" },
    @{ Row = 20; Col = 2; New = "923×5=4615" },
    @{ Row = 20; Col = 3; New = "274×6=1644" },
    @{ Row = 20; Col = 4; New = "516×6=3096" },
    @{ Row = 20; Col = 5; New = "244×6=1464" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
